# Generate Report for Handback
# Updates the timestamp strings recorded on the handback status report to
# reflect the new xliff generation / handoff / handback times.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
# "Latest HO Xliff Generate Date" for the a3b5a8b6-... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 03:04:16"

# --- zh-cn sheet ------------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the
# a3b5a8b6-... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 03:04:11"
$wsZhCn.Range("K2").Value = "2016-08-27 03:04:28"

# --- de-de sheet ------------------------------------------------------------
# "Correspond Handback DateTime" for the a3b5a8b6-... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-27 03:04:35"
